$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of test data: bestBatteryPhones library
$ws.Range("A10").Value = "bestBatteryPhones"
$ws.Range("B10").Value = "Infinix Hot 9 (Violet, 64 GB)"
$ws.Range("C10").Value = "Samsung Galaxy M11 (Black, 32 GB)"
$ws.Range("D10").Value = "Motorola G8 Power Lite (Royal Blue, 64 GB)"
$ws.Range("E10").Value = "Motorola Edge+ (Thunder Grey, 256 GB)"

# Apply the Menlo font style to B10 first, then fan the same formatting
# out to C10:E10 via a format-only paste so all four cells share a single
# new style entry (mirrors how Excel would coalesce identical formatting).
$fontCell = $ws.Range("B10")
$fontCell.Font.Name = "Menlo"
$fontCell.Font.Size = 11
$fontCell.Font.Color = 2236962

$fontCell.Copy()
$ws.Range("C10:E10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Widen columns C, D, E to fit the new content
$ws.Columns.Item(3).ColumnWidth = 39.333333333333336
$ws.Columns.Item(4).ColumnWidth = 50.5
$ws.Columns.Item(5).ColumnWidth = 41.166666666666664

# Select the new last cell, matching the authored selection
$ws.Range("E10").Select()
